$d = $word.ActiveDocument

# 1. "To run on localhost" -> "To run with CLI"
$d.Content.Find.Execute("To run on localhost", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "To run with CLI", 2) | Out-Null

# 2. "From logger directory" -> "From logger_publish directory"
$d.Content.Find.Execute("From logger directory, run command", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "From logger_publish directory, run command", 2) | Out-Null

# 3. "dotnet Logger.dll --urls" stays, but "From watcher directory" -> "From watcher_publish directory"
$d.Content.Find.Execute("From watcher directory, run command", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "From watcher_publish directory, run command", 2) | Out-Null

# 4. "Save file in watcher/watched directory" -> "Save file in watcher_publish/watched directory"
$d.Content.Find.Execute("Save file in watcher/watched directory", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Save file in watcher_publish/watched directory", 2) | Out-Null
